$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Professional summary: "all Black and Asian-American voters" -> "50M voters"
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2) | Out-Null

# -----------------------------------------------------------------
# 2) KEY ACHIEVEMENTS AND IMPACT -> "Impact" bullet list rewrite.
#    Locate the "KEY ACHIEVEMENTS AND IMPACT" heading, then its
#    "Impact" Heading3 sub-heading, then take the 4 bullet paragraphs
#    that directly follow it and replace them with the new 6-bullet
#    list (with bold spans).
# -----------------------------------------------------------------
$keyAchPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $keyAchPara = $i
        break
    }
}

$startPara = $keyAchPara + 2
$endPara = $startPara + 3

$rng = $d.Range($d.Paragraphs.Item($startPara).Range.Start, $d.Paragraphs.Item($endPara).Range.End)

$bullet1 = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **73.5%**"
$bullet2 = "• **`$4.7M** savings enabled nonprofit access"
$bullet3 = "• Legal precedent: Data analysis utilized in Supreme Court case"
$bullet4 = "• Expert methodology validated at highest judicial level"
$bullet5 = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"
$bullet6 = "• **178%** accuracy improvement in racial classification algorithms"

$fullText = $bullet1 + "`r" + $bullet2 + "`r" + $bullet3 + "`r" + $bullet4 + "`r" + $bullet5 + "`r" + $bullet6

$rngStart = $rng.Start
$rng.Text = $fullText

# Bold + color the three numeric spans, located by substring search
# relative to $rngStart so the offsets don't need to be hand-computed.
$boldTerms = @("73.5%", "`$4.7M", "178%")
foreach ($term in $boldTerms) {
    $offset = $fullText.IndexOf($term)
    $sub = $d.Range($rngStart + $offset, $rngStart + $offset + $term.Length)
    $sub.Font.Bold = 1
    $sub.Font.Color = 5258796
}

Write-Output "Done"
